# Documentation/AssetList.xlsx update:
# - Add reviewer (Steven) feedback notes in column F for a few rows
# - Mark the "ObjectiveComplete" row's Status as Completed (E15)
# These are entered in the same order the author typed them so the
# resulting shared-string table matches (notes added bottom-up, then
# the final "revised" note last).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# HeavyLanding (row 19) - Notes
$ws.Range("F19").Value = "Steven: Might be a bit too much of a funny sound"

# WeaponSwap (row 13) - Notes
$ws.Range("F13").Value = "Steven: A bit more bass"

# HoverbotMoving (row 5) - Notes
$ws.Range("F5").Value = "Steven: Might need higher frequency or layer (blends in too much)"

# ObjectiveComplete (row 15) - Status
$ws.Range("E15").Value = "Completed"

# Ambience (row 25) - Notes (objective-complete sound / lowered ambience volume revision)
$ws.Range("F25").Value = "Revised per Steven's feedback"

# Leave selection where the author finished editing
$ws.Range("F25").Select()
